$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 3 (R row)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 349
$wsOff.Range("C3").Value = 244
$wsOff.Range("D3").Value = 79
$wsOff.Range("E3").Value = 34
$wsOff.Range("F3").Value = 6

# Sheet "DEF" - row 3 (R row)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 379
$wsDef.Range("C3").Value = 256
$wsDef.Range("D3").Value = 102
$wsDef.Range("E3").Value = 51
